# Update error budget and detailed structural calculations to account for
# mispositioned bolt holes.
#
# - "structure" sheet: bracket material changes from continuous-glass-fiber
#   nylon to steel, and several geometry inputs for the cantilever-leg /
#   flange-leg calcs change (skin thickness, tensile strength, leg
#   dimensions, bolt-hole diameters). Downstream formulas recalc
#   automatically.
# - View/selection state moves from "joints and bearings" to "structure".

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("structure")
$ws2 = $wb.Worksheets.Item("joints and bearings")

# ---------------------------------------------------------------------
# structure sheet: data changes
# ---------------------------------------------------------------------

# Bracket material label: continuous glass fiber in nylon -> steel
$ws1.Range("D17").Value = "Material: steel"

# tensile strength (Mpa) for the new material
$ws1.Range("E19").Value = 370

# skin thickness (mm)
$ws1.Range("B20").Value = 1

# spacer thickness / B22 was "=B3-2*B20" -- bolt hole mispositioning means
# this is now an independent measured value instead of a derived one
$ws1.Range("B22").Value = 18

# B26 used to reference structure!B20/2; now derived from the updated
# spacer width (B22) and skin thickness (B20)
$ws1.Range("B26").Formula = "=(B22+2*B20)/2"

# sex bolt / flat head bolt hole diameters (cantilever leg - torsion calc)
$ws1.Range("B37").Value = 1
$ws1.Range("B38").Value = 45

# flat head bolt hole diameter (flange leg - axial loading calc)
$ws1.Range("B48").Value = 1

# ---------------------------------------------------------------------
# view / selection state: "structure" becomes the active sheet/tab,
# "joints and bearings" loses its tabSelected flag.
# ---------------------------------------------------------------------

$ws2.Activate()
$ws2.Range("A4").Select()
$ws2.Range("D13").Select()

$ws1.Activate()
$ws1.Range("A15").Select()
$ws1.Range("E30").Select()
